$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(11)
$shape.TextFrame.TextRange.Text = "Sexually transmitted in"
